# Applies the Turkish -> English localization edit described by the diff:
#   - rename the two worksheets
#   - translate the header row + risk-level labels + summary sheet labels
#   - move the active-cell selection on sheet 1
#   - re-size the header columns to the new (translated) best-fit widths

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Tahmin Sonuçları" -> "Forecast Results"
$ws2 = $wb.Worksheets.Item(2)   # "Risk Özeti"       -> "Risk Summary"

# ---- Sheet 1 : "Forecast Results" --------------------------------------

# Header row
$ws1.Range("A1").Value = "Country"
$ws1.Range("B1").Value = "Temperature_2013"
$ws1.Range("C1").Value = "Predicted_2100"
$ws1.Range("D1").Value = "Difference (°C)"
$ws1.Range("E1").Value = "Risk Level"
$ws1.Range("F1").Value = "Continent"

# Risk-level labels (column E), grouped by the contiguous blocks they occupy
$ws1.Range("E2:E31").Value = "Very High Risk"
$ws1.Range("E32:E38").Value = "High Risk"
$ws1.Range("E39:E57").Value = "Average Risk"
$ws1.Range("E58:E152").Value = "Low Risk"

# Column widths re-fitted for the (now longer/shorter) English headers
$ws1.Columns.Item(2).ColumnWidth = 17
$ws1.Columns.Item(3).ColumnWidth = 13.83
$ws1.Columns.Item(4).ColumnWidth = 13.5
$ws1.Columns.Item(5).ColumnWidth = 12.83

# Selection moved from F15 to I27
$ws1.Range("I27").Select()

# ---- Sheet 2 : "Risk Summary" ------------------------------------------

$ws2.Range("A1").Value = "Risk Level"
$ws2.Range("B1").Value = "Number of Countries"
$ws2.Range("A2").Value = "Low Risk"
$ws2.Range("A3").Value = "Very High Risk"
$ws2.Range("A4").Value = "Average Risk"
$ws2.Range("A5").Value = "High Risk"

# ---- Rename the sheets (do this last so the lookups above by index still
#      work regardless of whether the host resolves sheets by name) -------
$ws1.Name = "Forecast Results"
$ws2.Name = "Risk Summary"
